$d = $word.ActiveDocument

# The document currently ends with an otherwise-empty paragraph that only
# holds the "_GoBack" bookmark (right after the page break). We turn that
# paragraph into the centered "Script" heading and then add the body of
# the focus-group script as a series of new, indented paragraphs after it,
# keeping the bookmark on the very last paragraph.

$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range

$paragraphTexts = @(
    "Script",
    "Hello. Thank you for taking the time to be here today.",
    "Although we do not have a fully-working demonstration, we wanted to give you an idea of the planned features for Xpendit and to show you what we have for an interface design. ",
    "Xpendit is a group finance tracking app. Xpendit will allow you and your friends or roommates to keep track of what everyone owes each other, whether for a single night out or for an apartment.",
    "Xpendit’s differentiating features are the ability to have one-on-one transactions in addition to group transactions, as well as the ability to create shopping lists that can be edited by each member of a group. Not only will your roommates not forget to buy groceries again, but it will be easy to reimburse your friends so that nobody has to pay too much.",
    "Now that you’ve heard what Xpendit is all about, here is the app itself [show off the app].",
    "What do you all think of Xpendit’s feature set? Would you use the app yourself with your friends or for your apartment?",
    "What do you think about how Xpendit looks? ",
    "Thank you again for taking the time to answer our questions today."
)

$joined = [string]::Join("`r", $paragraphTexts)
$r.InsertBefore($joined)

$count = $d.Paragraphs.Count
$numNew = $paragraphTexts.Length
$firstNewIndex = $count - $numNew + 1

# First of the new paragraphs ("Script") becomes a centered heading.
$scriptPara = $d.Paragraphs.Item($firstNewIndex)
$scriptPara.Range.ParagraphFormat.Alignment = 1
$scriptPara.Range.Font.Size = 12

# Remaining new paragraphs get a first-line indent (0.5in = 36pt = 720 twips)
# and keep the same 24-half-point (12pt) run size used throughout the doc.
for ($i = $firstNewIndex + 1; $i -le $count; $i++) {
    $bodyPara = $d.Paragraphs.Item($i)
    $bodyPara.Range.ParagraphFormat.FirstLineIndent = 36
    $bodyPara.Range.Font.Size = 12
}
